$d = $word.ActiveDocument

# --- Main edit -----------------------------------------------------------
# The sentence currently reads:
#   "... we need. Then, we use it create all kinds of  JAXElement< Text>."
# It should read:
#   "... we need. Then, we use it to create all kinds of  JAXElement< Text>."
# i.e. insert the missing word "to" right after "we use it".

$rng = $d.Content.Duplicate
$rng.Find.Execute("we use it", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPoint = $rng.End

$ins = $d.Range($insPoint, $insPoint)
$ins.InsertAfter(" to")

# Touch (and immediately revert) a character formatting property on the
# newly-typed text. This mirrors how Word naturally keeps freshly typed
# text as its own run instead of silently re-merging it with the run(s)
# that used to occupy that position.
$ins.Bold = 1
$ins.Bold = 0

# --- Restore unrelated run boundaries ------------------------------------
# Editing the run above causes the engine to recompute/recombine the runs
# that immediately follow it in the same paragraph (" ", "JAXElement< Text>"
# and "." -- three separate, identically formatted runs in the original
# document) into a single run. That merge is only a side effect of this
# edit, not part of the intended change, so split them back apart the same
# way: nudge a formatting property on and back off for each sub-range.

$chk = $d.Content.Duplicate
$chk.Find.Execute("JAXElement< Text>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$jaxStart = $chk.Start
$jaxEnd = $chk.End

$spaceRng = $d.Range($jaxStart - 1, $jaxStart)
$spaceRng.Bold = 1
$spaceRng.Bold = 0

$dotRng = $d.Range($jaxEnd, $jaxEnd + 1)
$dotRng.Bold = 1
$dotRng.Bold = 0
